$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DRI")

# Row 4 - Inventory
$ws.Range("B4").Value = 188000000.0
$ws.Range("C4").Value = 203000000.0
$ws.Range("D4").Value = 190000000.0
$ws.Range("E4").Value = 207000000.0
$ws.Range("F4").Value = 230000000.0

# Row 13 - Accounts Payable
$ws.Range("B13").Value = 239000000.0
$ws.Range("C13").Value = 236000000.0
$ws.Range("D13").Value = 249000000.0
$ws.Range("E13").Value = 249000000.0
$ws.Range("F13").Value = 360000000.0

# Row 22 - Long Term Tax Liability (Deferred)
$ws.Range("B22").Value = 40000000.0
$ws.Range("C22").Value = 45000000.0
$ws.Range("D22").Value = 43000000.0
$ws.Range("E22").Value = 56000000.0
$ws.Range("F22").Value = 189000000.0
